$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B89").Value = "chayp"
$ws.Range("C89").Value = "bone"
$ws.Range("E89").Value = "noun"
$ws.Range("G89").Value = 1

$ws.Range("B90").Value = "ba_xbog_mgyemk"
$ws.Range("C90").Value = "butterfly"
$ws.Range("E90").Value = "noun"
$ws.Range("G90").Value = 1

$ws.Range("B91").Value = "wa_t'ukw"
$ws.Range("C91").Value = "lingcod"
$ws.Range("E91").Value = "noun"
$ws.Range("G91").Value = 1

$ws.Range("G90").Select() | Out-Null
